# Generate Report for Archive
# Two localization files (572e9e19-... and a01c5916-...) have moved from
# "Ready for handoff" into "In Translation" status. Update the per-language
# detail sheets (zh-cn, de-de) as well as the roll-up Overview sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: File Name / zh-cn / de-de columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus   # 572e9e19-...md, zh-cn
$overview.Range("C3").Value = $newStatus   # 572e9e19-...md, de-de
$overview.Range("B4").Value = $newStatus   # a01c5916-...md, zh-cn
$overview.Range("C4").Value = $newStatus   # a01c5916-...md, de-de

# --- zh-cn detail sheet: Status column (B) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus       # 572e9e19-...md
$zhcn.Range("B4").Value = $newStatus       # a01c5916-...md

# --- de-de detail sheet: Status column (B) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $newStatus       # 572e9e19-...md
$dede.Range("B4").Value = $newStatus       # a01c5916-...md
